# Auto-generated edit script applying scheduled market-price refresh to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 3494.6365
$ws.Cells.Item(51, 9).Value = 1780
$ws.Cells.Item(51, 10).Value = 3998.9412
$ws.Cells.Item(51, 11).Value = 1780
$ws.Cells.Item(51, 12).Value = 3998.9412
$ws.Cells.Item(51, 13).Value = -1296
$ws.Cells.Item(51, 14).Value = -4966.9412
$ws.Cells.Item(137, 8).Value = 1688.4138
$ws.Cells.Item(137, 9).Value = 1409.7142
$ws.Cells.Item(137, 10).Value = 2420
$ws.Cells.Item(137, 11).Value = 4229.142599999999
$ws.Cells.Item(137, 12).Value = 7260
$ws.Cells.Item(137, 13).Value = -1679.142599999999
$ws.Cells.Item(137, 14).Value = -12360
$ws.Cells.Item(138, 8).Value = 4880476
$ws.Cells.Item(138, 9).Value = 1293.3448
$ws.Cells.Item(138, 10).Value = 16671833
$ws.Cells.Item(138, 11).Value = 3880.0344
$ws.Cells.Item(138, 12).Value = 50015499
$ws.Cells.Item(138, 13).Value = 1259.9656
$ws.Cells.Item(138, 14).Value = -50025779
$ws.Cells.Item(141, 8).Value = 1071.25
$ws.Cells.Item(141, 9).Value = 1109.2858
$ws.Cells.Item(141, 10).Value = 805
$ws.Cells.Item(141, 11).Value = 3327.8574
$ws.Cells.Item(141, 12).Value = 2415
$ws.Cells.Item(141, 13).Value = 1852.1426
$ws.Cells.Item(141, 14).Value = -12775

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1026.1666
$ws.Cells.Item(61, 9).Value = 991.86487
$ws.Cells.Item(61, 10).Value = 1280
$ws.Cells.Item(61, 11).Value = 991.86487
$ws.Cells.Item(61, 12).Value = 1280
$ws.Cells.Item(61, 13).Value = -779.86487
$ws.Cells.Item(61, 14).Value = -1704
$ws.Cells.Item(74, 8).Value = 53504.715
$ws.Cells.Item(74, 9).Value = 59922.707
$ws.Cells.Item(74, 10).Value = 26228.25
$ws.Cells.Item(74, 11).Value = 59922.707
$ws.Cells.Item(74, 12).Value = 26228.25
$ws.Cells.Item(74, 13).Value = -59048.707
$ws.Cells.Item(74, 14).Value = -27976.25
$ws.Cells.Item(77, 8).Value = 53504.715
$ws.Cells.Item(77, 9).Value = 59922.707
$ws.Cells.Item(77, 10).Value = 26228.25
$ws.Cells.Item(77, 11).Value = 299613.535
$ws.Cells.Item(77, 12).Value = 131141.25
$ws.Cells.Item(77, 13).Value = -295245.535
$ws.Cells.Item(77, 14).Value = -139877.25
$ws.Cells.Item(88, 8).Value = 2433.3333
$ws.Cells.Item(88, 10).Value = 2300
$ws.Cells.Item(88, 12).Value = 2300
$ws.Cells.Item(88, 14).Value = -3112
$ws.Cells.Item(91, 8).Value = 2433.3333
$ws.Cells.Item(91, 10).Value = 2300
$ws.Cells.Item(91, 12).Value = 2300
$ws.Cells.Item(91, 13).Value = -1096
$ws.Cells.Item(91, 14).Value = -5108
$ws.Cells.Item(96, 8).Value = 25892
$ws.Cells.Item(96, 10).Value = 25892
$ws.Cells.Item(96, 12).Value = 25892
$ws.Cells.Item(96, 14).Value = -31384
$ws.Cells.Item(109, 8).Value = 30050
$ws.Cells.Item(109, 10).Value = 30050
$ws.Cells.Item(109, 12).Value = 30050
$ws.Cells.Item(109, 14).Value = -32824
$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 12).ClearContents()
$ws.Cells.Item(134, 14).Value = 0
$ws.Cells.Item(135, 8).Value = 30000
$ws.Cells.Item(135, 10).Value = 30000
$ws.Cells.Item(135, 12).Value = 30000
$ws.Cells.Item(135, 14).Value = -40140
$ws.Cells.Item(136, 8).Value = 1026.1666
$ws.Cells.Item(136, 9).Value = 991.86487
$ws.Cells.Item(136, 10).Value = 1280
$ws.Cells.Item(136, 11).Value = 2975.59461
$ws.Cells.Item(136, 12).Value = 3840
$ws.Cells.Item(136, 13).Value = -425.5946100000001
$ws.Cells.Item(136, 14).Value = -8940
$ws.Cells.Item(137, 8).Value = 39000
$ws.Cells.Item(137, 9).Value = 39000
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 39000
$ws.Cells.Item(137, 12).ClearContents()
$ws.Cells.Item(137, 14).Value = 0
$ws.Cells.Item(137, 13).Value = -33900
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 12).ClearContents()
$ws.Cells.Item(138, 14).Value = 0
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 12).ClearContents()
$ws.Cells.Item(139, 14).Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(57, 8).Value = 38666.668
$ws.Cells.Item(57, 10).Value = 38666.668
$ws.Cells.Item(57, 12).Value = 38666.668
$ws.Cells.Item(57, 14).Value = -40106.668
$ws.Cells.Item(99, 8).Value = 1701.25
$ws.Cells.Item(99, 9).Value = 1468.6666
$ws.Cells.Item(99, 11).Value = 1468.6666
$ws.Cells.Item(99, 13).Value = 29.33339999999998
$ws.Cells.Item(105, 8).Value = 2428.5715
$ws.Cells.Item(105, 9).Value = 2458.8235
$ws.Cells.Item(105, 10).Value = 2408
$ws.Cells.Item(105, 11).Value = 2458.8235
$ws.Cells.Item(105, 12).Value = 2408
$ws.Cells.Item(105, 13).Value = -711.8235
$ws.Cells.Item(105, 14).Value = -5902
$ws.Cells.Item(134, 8).Value = 4605.054
$ws.Cells.Item(134, 9).Value = 2936.682
$ws.Cells.Item(134, 10).Value = 7052
$ws.Cells.Item(134, 11).Value = 8810.045999999998
$ws.Cells.Item(134, 12).Value = 21156
$ws.Cells.Item(134, 13).Value = -6275.045999999998
$ws.Cells.Item(134, 14).Value = -26226
$ws.Cells.Item(136, 8).Value = 38666.668
$ws.Cells.Item(136, 10).Value = 38666.668
$ws.Cells.Item(136, 12).Value = 38666.668
$ws.Cells.Item(136, 14).Value = -48866.668
$ws.Cells.Item(137, 8).Value = 39939.395
$ws.Cells.Item(137, 10).Value = 39939.395
$ws.Cells.Item(137, 12).Value = 39939.395
$ws.Cells.Item(137, 14).Value = -50139.395
$ws.Cells.Item(138, 8).Value = 38770
$ws.Cells.Item(138, 10).Value = 38770
$ws.Cells.Item(138, 12).Value = 38770
$ws.Cells.Item(138, 14).Value = -49050
$ws.Cells.Item(139, 8).Value = 49214.285
$ws.Cells.Item(139, 10).Value = 48846.152
$ws.Cells.Item(139, 12).Value = 48846.152
$ws.Cells.Item(139, 14).Value = -59126.152

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 23257780
$ws.Cells.Item(31, 9).Value = 47619964
$ws.Cells.Item(31, 10).Value = 2970.0908
$ws.Cells.Item(31, 11).Value = 47619964
$ws.Cells.Item(31, 12).Value = 2970.0908
$ws.Cells.Item(31, 13).Value = -47619669
$ws.Cells.Item(31, 14).Value = -3560.0908
$ws.Cells.Item(34, 8).Value = 23257780
$ws.Cells.Item(34, 9).Value = 47619964
$ws.Cells.Item(34, 10).Value = 2970.0908
$ws.Cells.Item(34, 11).Value = 47619964
$ws.Cells.Item(34, 12).Value = 2970.0908
$ws.Cells.Item(34, 13).Value = -47619762
$ws.Cells.Item(34, 14).Value = -3374.0908
$ws.Cells.Item(58, 8).Value = 2678.9614
$ws.Cells.Item(58, 9).Value = 2882.739
$ws.Cells.Item(58, 10).Value = 1116.6666
$ws.Cells.Item(58, 11).Value = 2882.739
$ws.Cells.Item(58, 12).Value = 1116.6666
$ws.Cells.Item(58, 13).Value = -2679.739
$ws.Cells.Item(58, 14).Value = -1522.6666
$ws.Cells.Item(134, 8).Value = 33335670
$ws.Cells.Item(134, 9).Value = 4350378.5
$ws.Cells.Item(134, 10).Value = 100001840
$ws.Cells.Item(134, 11).Value = 13051135.5
$ws.Cells.Item(134, 12).Value = 300005520
$ws.Cells.Item(134, 13).Value = -13048600.5
$ws.Cells.Item(134, 14).Value = -300010590
$ws.Cells.Item(135, 8).Value = 34212.625
$ws.Cells.Item(135, 10).Value = 34212.625
$ws.Cells.Item(135, 12).Value = 34212.625
$ws.Cells.Item(135, 14).Value = -44352.625
$ws.Cells.Item(136, 8).Value = 2678.9614
$ws.Cells.Item(136, 9).Value = 2882.739
$ws.Cells.Item(136, 10).Value = 1116.6666
$ws.Cells.Item(136, 11).Value = 8648.217000000001
$ws.Cells.Item(136, 12).Value = 3349.9998
$ws.Cells.Item(136, 13).Value = -6098.217000000001
$ws.Cells.Item(136, 14).Value = -8449.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 906.35
$ws.Cells.Item(131, 10).Value = 908.04126
$ws.Cells.Item(131, 12).Value = 2724.12378
$ws.Cells.Item(131, 14).Value = -12804.12378

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 1798.7142
$ws.Cells.Item(126, 9).Value = 2261.7144
$ws.Cells.Item(126, 10).Value = 1335.7142
$ws.Cells.Item(126, 11).Value = 6785.1432
$ws.Cells.Item(126, 12).Value = 4007.1426
$ws.Cells.Item(126, 13).Value = -4315.1432
$ws.Cells.Item(126, 14).Value = -8947.142599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 1250
$ws.Cells.Item(100, 9).Value = 1000
$ws.Cells.Item(100, 10).Value = 1500
$ws.Cells.Item(100, 11).Value = 1000
$ws.Cells.Item(100, 12).Value = 1500
$ws.Cells.Item(100, 13).Value = -459
$ws.Cells.Item(100, 14).Value = -2582

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 897.6842
$ws.Cells.Item(100, 9).Value = 979.75
$ws.Cells.Item(100, 10).Value = 460
$ws.Cells.Item(100, 11).Value = 1959.5
$ws.Cells.Item(100, 12).Value = 920
$ws.Cells.Item(100, 13).Value = -1418.5
$ws.Cells.Item(100, 14).Value = -2002
